# Auto-generated Excel COM-interop script
# Applies scheduled market-price/profit refresh updates to Sheets/Zodiark_Profits.xlsx
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit* for leve-crafting rows)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 6046.4707
$ws.Range("I38").Value = 998
$ws.Range("J38").Value = 11726
$ws.Range("K38").Value = 2994
$ws.Range("L38").Value = 35178
$ws.Range("M38").Value = -2622
$ws.Range("N38").Value = -35922
$ws.Range("H39").Value = 742.82355
$ws.Range("I39").Value = 110.75
$ws.Range("J39").Value = 2259.8
$ws.Range("K39").Value = 332.25
$ws.Range("L39").Value = 6779.400000000001
$ws.Range("M39").Value = -36.25
$ws.Range("N39").Value = -7371.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2668.5
$ws.Range("J61").Value = 3500
$ws.Range("L61").Value = 3500
$ws.Range("N61").Value = -3924
$ws.Range("H136").Value = 2668.5
$ws.Range("J136").Value = 3500
$ws.Range("L136").Value = 10500
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1284.1666
$ws.Range("I64").Value = 996
$ws.Range("K64").Value = 996
$ws.Range("M64").Value = -771
$ws.Range("H67").Value = 1284.1666
$ws.Range("I67").Value = 996
$ws.Range("K67").Value = 996
$ws.Range("M67").Value = -216
$ws.Range("H86").Value = 2287.7917
$ws.Range("I86").Value = 1894.8334
$ws.Range("K86").Value = 1894.8334
$ws.Range("M86").Value = -771.8334
$ws.Range("H89").Value = 2287.7917
$ws.Range("I89").Value = 1894.8334
$ws.Range("K89").Value = 9474.167
$ws.Range("M89").Value = -3858.166999999999
$ws.Range("H99").Value = 102565260
$ws.Range("I99").Value = 133334320
$ws.Range("K99").Value = 133334320
$ws.Range("M99").Value = -133332822
$ws.Range("H105").Value = 1941.75
$ws.Range("I105").Value = 1910.1904
$ws.Range("K105").Value = 1910.1904
$ws.Range("M105").Value = -163.1904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2293.4736
$ws.Range("I31").Value = 1761.0625
$ws.Range("J31").Value = 5133
$ws.Range("K31").Value = 1761.0625
$ws.Range("L31").Value = 5133
$ws.Range("M31").Value = -1466.0625
$ws.Range("N31").Value = -5723
$ws.Range("H34").Value = 2293.4736
$ws.Range("I34").Value = 1761.0625
$ws.Range("J34").Value = 5133
$ws.Range("K34").Value = 1761.0625
$ws.Range("L34").Value = 5133
$ws.Range("M34").Value = -1559.0625
$ws.Range("N34").Value = -5537
$ws.Range("H58").Value = 6447
$ws.Range("I58").Value = 6878.25
$ws.Range("K58").Value = 6878.25
$ws.Range("M58").Value = -6675.25
$ws.Range("H105").Value = 9044.479
$ws.Range("I105").Value = 11098.889
$ws.Range("K105").Value = 11098.889
$ws.Range("M105").Value = -9351.889
$ws.Range("H107").Value = 598.6
$ws.Range("I107").Value = 585
$ws.Range("K107").Value = 585
$ws.Range("M107").Value = 1335
$ws.Range("H134").Value = 1397.6666
$ws.Range("I134").Value = 1138.2
$ws.Range("J134").Value = 1722
$ws.Range("K134").Value = 3414.6
$ws.Range("L134").Value = 5166
$ws.Range("M134").Value = -879.6000000000004
$ws.Range("N134").Value = -10236
$ws.Range("H136").Value = 6447
$ws.Range("I136").Value = 6878.25
$ws.Range("K136").Value = 20634.75
$ws.Range("M136").Value = -18084.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 243.22223
$ws.Range("I5").Value = 243.22223
$ws.Range("K5").Value = 729.66669
$ws.Range("M5").Value = -617.66669
$ws.Range("H10").Value = 133.6
$ws.Range("I10").Value = 133.6
$ws.Range("K10").Value = 400.8
$ws.Range("M10").Value = -261.8
$ws.Range("H21").Value = 998.1667
$ws.Range("I21").Value = 622.5
$ws.Range("J21").Value = 1749.5
$ws.Range("K21").Value = 1867.5
$ws.Range("L21").Value = 5248.5
$ws.Range("M21").Value = -1694.5
$ws.Range("N21").Value = -5594.5
$ws.Range("H131").Value = 515
$ws.Range("I131").Value = 515
$ws.Range("K131").Value = 1545
$ws.Range("M131").Value = 3495
$ws.Range("H135").Value = 243.22223
$ws.Range("I135").Value = 243.22223
$ws.Range("K135").Value = 2189.00007
$ws.Range("M135").Value = 345.9999299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 999.7143
$ws.Range("I97").Value = 774.5
$ws.Range("J97").Value = 1089.8
$ws.Range("K97").Value = 774.5
$ws.Range("L97").Value = 1089.8
$ws.Range("M97").Value = -278.5
$ws.Range("N97").Value = -2081.8
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 52637180
$ws.Range("I7").Value = 83337640
$ws.Range("J7").Value = 7827.7144
$ws.Range("K7").Value = 83337640
$ws.Range("L7").Value = 7827.7144
$ws.Range("M7").Value = -83337528
$ws.Range("N7").Value = -8051.7144
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 10000
$ws.Range("K13").Value = 10000
$ws.Range("M13").Value = -9860
$ws.Range("H82").Value = 32260378
$ws.Range("I82").Value = 52633236
$ws.Range("K82").Value = 52633236
$ws.Range("M82").Value = -52632875
$ws.Range("H85").Value = 32260378
$ws.Range("I85").Value = 52633236
$ws.Range("K85").Value = 52633236
$ws.Range("M85").Value = -52631988
$ws.Range("H122").Value = 11653.272
$ws.Range("I122").Value = 5183.143
$ws.Range("J122").Value = 22976
$ws.Range("K122").Value = 15549.429
$ws.Range("L122").Value = 68928
$ws.Range("M122").Value = -13099.429
$ws.Range("N122").Value = -73828
$ws.Range("H126").Value = 52637180
$ws.Range("I126").Value = 83337640
$ws.Range("J126").Value = 7827.7144
$ws.Range("K126").Value = 250012920
$ws.Range("L126").Value = 23483.1432
$ws.Range("M126").Value = -250010450
$ws.Range("N126").Value = -28423.1432
$ws.Range("H136").Value = 4065.6667
$ws.Range("I136").Value = 3799.6667
$ws.Range("J136").Value = 4331.6665
$ws.Range("K136").Value = 11399.0001
$ws.Range("L136").Value = 12994.9995
$ws.Range("M136").Value = -8849.000100000001
$ws.Range("N136").Value = -18094.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8163.1177
$ws.Range("I136").Value = 9236.462
$ws.Range("K136").Value = 27709.386
$ws.Range("M136").Value = -25159.386
